# "add until User Interface" -- rework Requirements sheet: expand the
# single "User Interface" section into the full set of sections
# (SystemAdministor / Course Administrator / Course Staff / All Users /
# Workflow Management Activities / User Management / Exam Management)
# spread out over rows 1-58, with blank spacer rows between sections.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$FILL = 49407      # RGB(255,192,0) == fgColor FFFFC000 used by the existing section headers
$VCENTER = -4108   # xlCenter (vertical)
$VBOTTOM = -4107   # xlBottom (default / "general")

# ---------------------------------------------------------------------
# 0. Drop the three old content rows that no longer exist anywhere in
#    the new layout (old rows 4 & 5 collapse away) and the stray B:D
#    cells of what used to be section-header rows 8/11/12 (they turn
#    into plain single-column rows in the new layout).
# ---------------------------------------------------------------------
$ws.Range("A4:A5").ClearContents()
$ws.Range("A4:A5").Style = "Normal"
$ws.Range("B8:D8").ClearContents()
$ws.Range("B8:D8").Style = "Normal"
$ws.Range("B11:D12").ClearContents()
$ws.Range("B11:D12").Style = "Normal"

# ---------------------------------------------------------------------
# 1. Column layout: column A keeps its width but loses the leftover
#    per-column style; only column E keeps an explicit (default) width
#    now, instead of B:XFD.
# ---------------------------------------------------------------------
$ws.Columns("A:A").ClearFormats()
$ws.Columns("B:XFD").ClearFormats()
$ws.Columns.Item(5).ColumnWidth = 8.88671875

# ---------------------------------------------------------------------
# 2. Row 1 (header) is untouched except that D1 keeps the same text
#    ("Priority") -- nothing to change there content-wise.
# ---------------------------------------------------------------------

# ---------------------------------------------------------------------
# 3. Row 2 "User Interface:" section header -- already existed, just
#    re-assert the fill/valign so it matches after the column-format
#    reset above.
# ---------------------------------------------------------------------
function Set-SectionHeader($rng, $vAlign) {
    $rng.Interior.Color = $FILL
    $rng.VerticalAlignment = $vAlign
}

function Clear-PlainCell($rng) {
    $rng.Interior.ColorIndex = -4142
    $rng.VerticalAlignment = -4107
}

$ws.Range("A2:D2").Value = "User Interface:"
Set-SectionHeader $ws.Range("A2:D2") $VCENTER
$ws.Range("B2").Value = ""
$ws.Range("C2").Value = ""
$ws.Range("D2").Value = ""
$ws.Range("A2").Value = "User Interface:"

# Row 3: single requirement line under "User Interface:"
$ws.Range("A3").Value = "Web-based"
Clear-PlainCell $ws.Range("A3")

# ---------------------------------------------------------------------
# 4. SystemAdministor section (rows 6-13)
# ---------------------------------------------------------------------
$ws.Range("A6").Value = "SystemAdministor"
$ws.Range("B6").Value = ""
$ws.Range("C6").Value = ""
$ws.Range("D6").Value = ""
Set-SectionHeader $ws.Range("A6:D6") $VCENTER

$sysAdminRows = @(
    @(7,  "Manage WMS at the dept/university level"),
    @(8,  "Migrate, install, update, manage WMS system"),
    @(9,  "Manage courses in WMS"),
    @(10, "Create new course entry in WMS"),
    @(11, "Backup and remove course entry in WMS"),
    @(12, "Assign roles to courses"),
    @(13, "Remove, change roles per course (when course staff changes)")
)
foreach ($r in $sysAdminRows) {
    $cell = $ws.Cells.Item($r[0], 1)
    $cell.Value = $r[1]
    Clear-PlainCell $cell
}

# ---------------------------------------------------------------------
# 5. Course Administrator section (rows 16-22)
# ---------------------------------------------------------------------
$ws.Range("A16").Value = "Course Administrator"
$ws.Range("B16").Value = ""
$ws.Range("C16").Value = ""
$ws.Range("D16").Value = ""
Set-SectionHeader $ws.Range("A16:D16") $VCENTER

$courseAdminRows = @(
    @(17, "Create tasks"),
    @(18, "Define course staff"),
    @(19, "Assign roles to course staff"),
    @(20, "Define exam"),
    @(21, "Generate exam documents"),
    @(22, "Inspect changes by course staff")
)
foreach ($r in $courseAdminRows) {
    $cell = $ws.Cells.Item($r[0], 1)
    $cell.Value = $r[1]
    Clear-PlainCell $cell
}

# ---------------------------------------------------------------------
# 6. Course Staff section (rows 25-28)
# ---------------------------------------------------------------------
$ws.Range("A25").Value = "Course Staff"
$ws.Range("B25").Value = ""
$ws.Range("C25").Value = ""
$ws.Range("D25").Value = ""
Set-SectionHeader $ws.Range("A25:D25") $VCENTER

$courseStaffRows = @(
    @(26, "Roles (Junior, Associate, Senior)"),
    @(27, "Activities"),
    @(28, "People (Instructores TAs  Graders)")
)
foreach ($r in $courseStaffRows) {
    $cell = $ws.Cells.Item($r[0], 1)
    $cell.Value = $r[1]
    Clear-PlainCell $cell
}

# ---------------------------------------------------------------------
# 7. All Users section (rows 30-34)
# ---------------------------------------------------------------------
$ws.Range("A30").Value = "All Users"
$ws.Range("B30").Value = ""
$ws.Range("C30").Value = ""
$ws.Range("D30").Value = ""
Set-SectionHeader $ws.Range("A30:D30") $VCENTER

$allUsersRows = @(
    @(31, "Watch pending tasks"),
    @(32, "Ask for work"),
    @(33, "Perform specific task"),
    @(34, "Ask for a task")
)
foreach ($r in $allUsersRows) {
    $cell = $ws.Cells.Item($r[0], 1)
    $cell.Value = $r[1]
    Clear-PlainCell $cell
}

# ---------------------------------------------------------------------
# 8. Workflow Management Activities section (rows 44-47)
#    -- this header keeps the "old" fill-only style (no forced vertical
#    centering), same as the pre-existing Exam/User Management headers.
# ---------------------------------------------------------------------
$ws.Range("A44").Value = "Workflow Management Activities:"
$ws.Range("B44").Value = ""
$ws.Range("C44").Value = ""
$ws.Range("D44").Value = ""
Set-SectionHeader $ws.Range("A44:D44") $VBOTTOM

$workflowRows = @(
    @(45, "Log into system"),
    @(46, "Maintain history (continue where left off)"),
    @(47, "What controls urgency")
)
foreach ($r in $workflowRows) {
    $cell = $ws.Cells.Item($r[0], 1)
    $cell.Value = $r[1]
    Clear-PlainCell $cell
}

# ---------------------------------------------------------------------
# 9. User Management section (rows 55-57) -- pre-existing section,
#    moved further down the sheet.
# ---------------------------------------------------------------------
$ws.Range("A55").Value = "User Management:"
$ws.Range("B55").Value = ""
$ws.Range("C55").Value = ""
$ws.Range("D55").Value = ""
Set-SectionHeader $ws.Range("A55:D55") $VBOTTOM

$userMgmtRows = @(
    @(56, "Username-Password security"),
    @(57, "Roles")
)
foreach ($r in $userMgmtRows) {
    $cell = $ws.Cells.Item($r[0], 1)
    $cell.Value = $r[1]
    Clear-PlainCell $cell
}

# ---------------------------------------------------------------------
# 10. Exam Management section header (row 58) -- pre-existing section,
#     now the very last row of the sheet.
# ---------------------------------------------------------------------
$ws.Range("A58").Value = "Exam Management:"
$ws.Range("B58").Value = ""
$ws.Range("C58").Value = ""
$ws.Range("D58").Value = ""
Set-SectionHeader $ws.Range("A58:D58") $VBOTTOM

# ---------------------------------------------------------------------
# 11. View state: scroll/selection now sits on the last section header.
# ---------------------------------------------------------------------
$ws.Range("A55:D58").Select()

